$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "data refreshed at" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 11:22"

# --- Refresh daily COVID numbers for countries whose ranking did not change ---

# Alemania (row 25)
$ws.Range("B25").Value = 273523
$ws.Range("C25").Value = 46
$ws.Range("E25").Value = 20053

# Indonesia (row 26)
$ws.Range("B26").Value = 248852
$ws.Range("C26").Value = 4176
$ws.Range("D26").Value = 180797
$ws.Range("E26").Value = 58378
$ws.Range("G26").Value = 124
$ws.Range("H26").Value = 9677

# Ucrania (row 28)
$ws.Range("D28").Value = 78184
$ws.Range("E28").Value = 96586

# Polonia (row 47)
$ws.Range("B47").Value = 79988
$ws.Range("C47").Value = 748
$ws.Range("D47").Value = 64604
$ws.Range("E47").Value = 13086
$ws.Range("G47").Value = 5
$ws.Range("H47").Value = 2298

# Austria (row 69)
$ws.Range("B69").Value = 38658
$ws.Range("C69").Value = 563
$ws.Range("D69").Value = 29516
$ws.Range("E69").Value = 8375
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 767

# Croacia (row 89)
$ws.Range("B89").Value = 14992
$ws.Range("C89").Value = 70
$ws.Range("D89").Value = 12737
$ws.Range("E89").Value = 2002
$ws.Range("G89").Value = 5
$ws.Range("H89").Value = 253

# Finlandia (row 103)
$ws.Range("B103").Value = 9046
$ws.Range("C103").Value = 66
$ws.Range("E103").Value = 1007

# Eslovaquia (row 111)
$ws.Range("B111").Value = 6756
$ws.Range("C111").Value = 79
$ws.Range("D111").Value = 3571
$ws.Range("E111").Value = 3146

# Hong Kong (row 119)
$ws.Range("B119").Value = 5039
$ws.Range("C119").Value = 6
$ws.Range("D119").Value = 4717
$ws.Range("E119").Value = 219

# Eslovenia (row 128)
$ws.Range("B128").Value = 4470
$ws.Range("C128").Value = 50
$ws.Range("D128").Value = 3048
$ws.Range("E128").Value = 1280

# --- Lituania overtakes Siria (rows 131/132 swap order, Lituania gets new numbers) ---
$ws.Range("A131").Value = "Lituania"
$ws.Range("B131").Value = 3814
$ws.Range("C131").Value = 70
$ws.Range("D131").Value = 2199
$ws.Range("E131").Value = 1528
$ws.Range("H131").Value = 87

$ws.Range("A132").Value = "Siria"
$ws.Range("B132").Value = 3800
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 946
$ws.Range("E132").Value = 2682
$ws.Range("H132").Value = 172

# --- Montserrat overtakes Islas Malvinas (rows 214/215 swap order) ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
